# Extend the results table (vm_pu) from 13 columns (0..12) to 17 columns (0..16):
# add header values 13,14,15,16 in O1:R1 (bold/border/center style, like existing headers),
# refresh the recomputed per-bus voltage magnitudes in columns C,D,E,I,J,K,L,N (rows 2-25),
# and add two new result columns Q (constant 1.02) and R (new computed values) for rows 2-25.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- New header cells O1:P1:Q1:R1, matching the style of the existing header cells (e.g. N1) ---
$ws.Range("N1").Copy()
$ws.Range("O1:R1").PasteSpecial(-4122)  # xlPasteFormats
$excel.CutCopyMode = $false
$ws.Range("O1").Value = 13
$ws.Range("P1").Value = 14
$ws.Range("Q1").Value = 15
$ws.Range("R1").Value = 16

$cols = @("C","D","E","I","J","K","L","N","Q","R")
$data = @(
    @(0.9993172784004751,1.021153796773663,1.005371000175455,1.047252258917627,1.021518428746962,1.032326086766614,1.016756012165627,1.022969102085584,1.02,1.033927555859995),
    @(1.005739161254685,1.025647163740762,1.010683626841256,1.049039085722464,1.02607458434123,1.035962073072042,1.021181609590869,1.027531727943394,1.02,1.03649587251796),
    @(1.009783940206496,1.028479223287288,1.014036239771889,1.050149006701047,1.028940640336147,1.038245662753256,1.023968160714857,1.030401854066565,1.02,1.038111482905451),
    @(1.011468692378906,1.029661122839955,1.01543721107426,1.050609278103296,1.03013604399037,1.039198307136427,1.02513280508566,1.031598955331092,1.02,1.038792271050643),
    @(1.011761671467029,1.029868980408191,1.015684600655205,1.050691235185495,1.03034694855333,1.03936763925606,1.025340465735508,1.031810159402722,1.02,1.038920761001531),
    @(1.00983788456975,1.028523193721598,1.014091192631628,1.050169495119014,1.028987240820612,1.038286248011975,1.024019527520492,1.030448520729064,1.02,1.038160335855932),
    @(1.00155164729018,1.022724126302878,1.007231722483096,1.047884051824779,1.0231149351058,1.033604837159864,1.018314618503826,1.024567875666574,1.02,1.034854783342521),
    @(0.9859867599324414,1.011842300709696,0.9943849166424276,1.043460264568035,1.012044254976849,1.024746122896897,1.007569756362888,1.013481473902015,1.02,1.028587749478976),
    @(0.9749099761258737,1.004122096424771,0.985291963390957,1.040236107188704,1.004155673473599,1.01841732187265,0.9999319527735625,1.005581689708197,1.02,1.024129783017159),
    @(0.9699659643628109,1.000689918986226,0.9812634095941459,1.038785878988771,1.000643861887691,1.01559883616365,0.9965475372486526,1.002064890947027,1.02,1.022170304430479),
    @(0.9680730685520957,0.9993733004428675,0.9797133995736602,1.03822477480804,0.9992929096749885,1.014512516884108,0.9952398598749699,1.000712020227226,1.02,1.021402215496618),
    @(0.9684755643217298,0.9996524777262582,0.9800409793482133,1.038343555051712,0.9995789716695013,1.014742282442272,0.9955154619707512,1.000998488462582,1.02,1.021562134297352),
    @(0.9698080309949866,1.000579730672563,0.9811332325012941,1.038738839381524,1.00053063593098,1.015507676917636,0.9964373894403536,1.00195150419647,1.02,1.02210476279498),
    @(0.9706345560823938,1.00115651957853,0.9818147568025686,1.038984981918119,1.001123301065386,1.015984852028369,0.9970140798873185,1.002545010983347,1.02,1.022448127174121),
    @(0.9753175316331362,1.004417990925721,0.9856562454484868,1.040366610072993,1.004465812942503,1.01867148805129,1.000251946555927,1.00589226961072,1.02,1.024350936026014),
    @(0.9782006204318076,1.006427921411466,0.9880243334863109,1.041212738326235,1.006523082702683,1.020324176810739,1.002245649084845,1.007952460929933,1.02,1.02552209177624),
    @(0.9798383353359409,1.00756637723276,0.989361208247312,1.041688251793185,1.007685276656525,1.021255770885233,1.003366143256993,1.009116305332494,1.02,1.026168969851055),
    @(0.9804130865646619,1.007969087618906,0.9898378863127615,1.041857592541265,1.008097957216106,1.021587854336732,1.003768868044628,1.009529571945805,1.02,1.026410249620082),
    @(0.9778912912778389,1.006211901366657,0.9877694131185328,1.041121856998303,1.006301897591135,1.020146394401867,1.002030786808624,1.007730961710155,1.02,1.025395036585729),
    @(0.9694446039521004,1.000330842215613,0.980845948284661,1.038634682911837,1.000277810792031,1.01530594939978,0.9961995190481875,1.001698320016825,1.02,1.021975884526061),
    @(0.9639116135522388,0.9964808110073449,0.9763080893830971,1.036986301546006,0.9963217441564162,1.012122307577474,0.9923639655469564,0.9977366353128194,1.02,1.019711106781021),
    @(0.9668295324776265,0.9985055888664306,0.9786872085580048,1.037852668591987,0.9984000507548166,1.013793209160936,0.9943701233188663,0.9998178933449751,1.02,1.02088272554292),
    @(0.9779777129824803,1.006264093397837,0.987820836485611,1.041139948205881,1.006350664423142,1.020182090348034,1.0020652617963,1.007779797796658,1.02,1.02539248024421),
    @(0.9901819397747642,1.014781038785437,0.9978594631919387,1.044674533493903,1.015044299457299,1.027154150982444,1.010491099368449,1.016485778789737,1.02,1.030319287274519)
)

for ($i = 0; $i -lt $data.Count; $i++) {
    $row = $i + 2
    for ($j = 0; $j -lt $cols.Count; $j++) {
        $ws.Range("$($cols[$j])$row").Value = $data[$i][$j]
    }
}

